# Auto-generated edit script applying the scheduled-runner Sheets update.
# For each affected row (identified by sheet + row number), update the
# pricing/profit columns (H-N) to the new computed values. Some rows also
# gain or lose a cell in columns M/N (profit fields) as part of the update;
# those are modeled as ClearContents() (remove) or a fresh .Value assignment
# (add) on the previously-blank cell.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 172966.58
$ws.Range("I6").Value = 201769.33
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 605307.99
$ws.Range("L6").Value = 450
$ws.Range("M6").Value = -605195.99
$ws.Range("N6").Value = -674
# Row 49
$ws.Range("H49").Value = 232.33333
$ws.Range("I49").Value = 232.33333
$ws.Range("K49").Value = 696.99999
$ws.Range("M49").Value = -560.99999
# Row 64
$ws.Range("H64").Value = 7634.727
$ws.Range("I64").Value = 7830.3335
$ws.Range("K64").Value = 7830.3335
$ws.Range("M64").Value = -7582.3335
# Row 67
$ws.Range("H67").Value = 7634.727
$ws.Range("I67").Value = 7830.3335
$ws.Range("K67").Value = 7830.3335
$ws.Range("M67").Value = -6972.3335
# Row 138
$ws.Range("H138").Value = 1289.5217
$ws.Range("J138").Value = 2997.5
$ws.Range("L138").Value = 8992.5
$ws.Range("N138").Value = -19272.5
# Row 141
$ws.Range("H141").Value = 7555.4
$ws.Range("I141").Value = 8194.583000000001
$ws.Range("K141").Value = 24583.749
$ws.Range("M141").Value = -19403.749

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3477.0908
$ws.Range("I45").Value = 2992.5715
$ws.Range("J45").Value = 4325
$ws.Range("K45").Value = 2992.5715
$ws.Range("L45").Value = 4325
$ws.Range("M45").Value = -2615.5715
$ws.Range("N45").Value = -5079
# Row 61
$ws.Range("H61").Value = 3100.6365
$ws.Range("I61").Value = 2831.3845
$ws.Range("K61").Value = 2831.3845
$ws.Range("M61").Value = -2619.3845
# Row 63
$ws.Range("H63").Value = 1794.2858
$ws.Range("I63").Value = 1863.5454
$ws.Range("J63").Value = 1718.1
$ws.Range("K63").Value = 1863.5454
$ws.Range("L63").Value = 1718.1
$ws.Range("M63").Value = -1177.5454
$ws.Range("N63").Value = -3090.1
# Row 66
$ws.Range("H66").Value = 1794.2858
$ws.Range("I66").Value = 1863.5454
$ws.Range("J66").Value = 1718.1
$ws.Range("K66").Value = 9317.726999999999
$ws.Range("L66").Value = 8590.5
$ws.Range("M66").Value = -5885.726999999999
$ws.Range("N66").Value = -15454.5
# Row 97
$ws.Range("H97").Value = 43480140
$ws.Range("I97").Value = 1361.2
$ws.Range("J97").Value = 333338660
$ws.Range("K97").Value = 1361.2
$ws.Range("L97").Value = 333338660
$ws.Range("M97").Value = -865.2
$ws.Range("N97").Value = -333339652
# Row 102
$ws.Range("H102").Value = 3787.125
$ws.Range("I102").Value = 3787.125
$ws.Range("K102").Value = 3787.125
$ws.Range("M102").Value = -2165.125
# Row 122
$ws.Range("H122").Value = 2345.625
$ws.Range("I122").Value = 1910.4
$ws.Range("K122").Value = 5731.200000000001
$ws.Range("M122").Value = -3281.200000000001
# Row 132
$ws.Range("H132").Value = 464320.75
$ws.Range("I132").Value = 482010.03
$ws.Range("K132").Value = 1446030.09
$ws.Range("M132").Value = -1443500.09
# Row 136
$ws.Range("H136").Value = 3100.6365
$ws.Range("I136").Value = 2831.3845
$ws.Range("K136").Value = 8494.1535
$ws.Range("M136").Value = -5944.1535

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 200.33333
$ws.Range("J22").Value = 200
$ws.Range("L22").Value = 200
$ws.Range("N22").Value = -546
# Row 134
$ws.Range("H134").Value = 5012.769
$ws.Range("I134").Value = 2403.1614
$ws.Range("K134").Value = 7209.4842
$ws.Range("M134").Value = -4674.4842

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2459.3333
$ws.Range("I31").Value = 2382.0952
$ws.Range("K31").Value = 2382.0952
$ws.Range("M31").Value = -2087.0952
# Row 34
$ws.Range("H34").Value = 2459.3333
$ws.Range("I34").Value = 2382.0952
$ws.Range("K34").Value = 2382.0952
$ws.Range("M34").Value = -2180.0952
# Row 86
$ws.Range("H86").Value = 22614.285
$ws.Range("J86").Value = 6150
$ws.Range("L86").Value = 6150
$ws.Range("N86").Value = -8396
# Row 89
$ws.Range("H89").Value = 22614.285
$ws.Range("J89").Value = 6150
$ws.Range("L89").Value = 30750
$ws.Range("N89").Value = -41982
# Row 105
$ws.Range("H105").Value = 23515.2
$ws.Range("I105").Value = 33859.668
$ws.Range("K105").Value = 33859.668
$ws.Range("M105").Value = -32112.668
# Row 132
$ws.Range("H132").Value = 3447.7778
$ws.Range("I132").Value = 3360.4285
$ws.Range("J132").Value = 3753.5
$ws.Range("K132").Value = 10081.2855
$ws.Range("L132").Value = 11260.5
$ws.Range("M132").Value = -7551.2855
$ws.Range("N132").Value = -16320.5
# Row 134
$ws.Range("H134").Value = 1973.5834
$ws.Range("I134").Value = 1507.9524
$ws.Range("J134").Value = 5233
$ws.Range("K134").Value = 4523.857199999999
$ws.Range("L134").Value = 15699
$ws.Range("M134").Value = -1988.857199999999
$ws.Range("N134").Value = -20769

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 94.625
$ws.Range("I2").Value = 107
$ws.Range("K2").Value = 642
$ws.Range("M2").Value = -529
# Row 105
$ws.Range("H105").Value = 21142.857
$ws.Range("I105").Value = 15000
$ws.Range("K105").Value = 45000
$ws.Range("M105").Value = -42379
# Row 136
$ws.Range("H136").Value = 12915.333
$ws.Range("I136").Value = 13000
$ws.Range("K136").Value = 39000
$ws.Range("M136").Value = -33900

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 4762024.5
$ws.Range("I2").Value = 8333382.5
$ws.Range("J2").Value = 213.88889
$ws.Range("K2").Value = 8333382.5
$ws.Range("L2").Value = 213.88889
$ws.Range("M2").Value = -8333269.5
$ws.Range("N2").Value = -439.88889
# Row 20
$ws.Range("H20").Value = 17633
$ws.Range("J20").Value = 11449.5
$ws.Range("L20").Value = 11449.5
$ws.Range("N20").Value = -11939.5
# Row 21
$ws.Range("H21").Value = 10000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 10000
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -10346
# Row 24
$ws.Range("H24").Value = 10666.667
$ws.Range("J24").Value = 10666.667
$ws.Range("L24").Value = 10666.667
$ws.Range("N24").Value = -11012.667
# Row 30
$ws.Range("H30").Value = 10000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 10000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 10000
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -10210
# Row 102
$ws.Range("H102").Value = 2655
$ws.Range("I102").Value = 2655
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2655
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1033
$ws.Range("N102").ClearContents()
# Row 122
$ws.Range("H122").Value = 46725.434
$ws.Range("J122").Value = 4833
$ws.Range("L122").Value = 14499
$ws.Range("N122").Value = -19399
# Row 126
$ws.Range("H126").Value = 2974.375
$ws.Range("J126").Value = 3099.75
$ws.Range("L126").Value = 9299.25
$ws.Range("N126").Value = -14239.25

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3700
$ws.Range("I7").Value = 3700
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3700
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3588
$ws.Range("N7").ClearContents()
# Row 40
$ws.Range("H40").Value = 4128.222
$ws.Range("I40").Value = 3632.8
$ws.Range("J40").Value = 4747.5
$ws.Range("K40").Value = 3632.8
$ws.Range("L40").Value = 4747.5
$ws.Range("M40").Value = -3496.8
$ws.Range("N40").Value = -5019.5
# Row 68
$ws.Range("H68").Value = 1861.4615
$ws.Range("I68").Value = 1875
$ws.Range("J68").Value = 1699
$ws.Range("K68").Value = 1875
$ws.Range("L68").Value = 1699
$ws.Range("M68").Value = -1126
$ws.Range("N68").Value = -3197
# Row 71
$ws.Range("H71").Value = 1861.4615
$ws.Range("I71").Value = 1875
$ws.Range("J71").Value = 1699
$ws.Range("K71").Value = 9375
$ws.Range("L71").Value = 8495
$ws.Range("M71").Value = -5631
$ws.Range("N71").Value = -15983
# Row 122
$ws.Range("H122").Value = 7315.6665
$ws.Range("I122").Value = 4949.5
$ws.Range("J122").Value = 8498.75
$ws.Range("K122").Value = 14848.5
$ws.Range("L122").Value = 25496.25
$ws.Range("M122").Value = -12398.5
$ws.Range("N122").Value = -30396.25
# Row 126
$ws.Range("H126").Value = 3700
$ws.Range("I126").Value = 3700
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11100
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8630
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 5002829.5
$ws.Range("I132").Value = 25001596
$ws.Range("J132").Value = 3137.5
$ws.Range("K132").Value = 75004788
$ws.Range("L132").Value = 9412.5
$ws.Range("M132").Value = -75002258
$ws.Range("N132").Value = -14472.5
# Row 136
$ws.Range("H136").Value = 14631.444
$ws.Range("I136").Value = 5986.5
$ws.Range("J136").Value = 17101.428
$ws.Range("K136").Value = 17959.5
$ws.Range("L136").Value = 51304.284
$ws.Range("M136").Value = -15409.5
$ws.Range("N136").Value = -56404.284
# Row 139
$ws.Range("H139").Value = 89999.8
$ws.Range("J139").Value = 89999.8
$ws.Range("L139").Value = 89999.8
$ws.Range("N139").Value = -100279.8

$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Range("H40").Value = 2666.6667
$ws.Range("I40").Value = 2666.6667
$ws.Range("K40").Value = 2666.6667
$ws.Range("M40").Value = -2517.6667
# Row 122
$ws.Range("H122").Value = 46904.48
$ws.Range("I122").Value = 1794.1177
$ws.Range("K122").Value = 5382.3531
$ws.Range("M122").Value = -2932.3531
# Row 126
$ws.Range("H126").Value = 2195.3333
$ws.Range("I126").Value = 2195.3333
$ws.Range("K126").Value = 6585.999899999999
$ws.Range("M126").Value = -4115.999899999999
# Row 132
$ws.Range("H132").Value = 5998.75
$ws.Range("I132").Value = 5998.75
$ws.Range("K132").Value = 17996.25
$ws.Range("M132").Value = -15466.25
# Row 136
$ws.Range("H136").Value = 1343.8462
$ws.Range("I136").Value = 1345.3043
$ws.Range("J136").Value = 1332.6666
$ws.Range("K136").Value = 4035.9129
$ws.Range("L136").Value = 3997.9998
$ws.Range("M136").Value = -1485.9129
$ws.Range("N136").Value = -9097.9998
